$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 5 corresponds to the
# 968b3c4b-617e-460b-ac64-45d96ecf4a67 handback entry.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-03-02 15:06:11"
$wsZh.Range("G5").Value = "2016-03-02 15:07:07"

# de-de sheet: same row/entry, different locale.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-03-02 15:06:26"
$wsDe.Range("G5").Value = "2016-03-02 15:07:30"
